# Add the "ItemValue" lookup sheet (after "Crystals") and wire up the
# Crystals sheet with a "sacrifice conversion" table (rows 18-23) that
# looks values up from it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "ItemValue", inserted as the last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$iv = $wb.Worksheets.Add($null, $lastSheet)
$iv.Name = "ItemValue"

# Fill it in roughly in the order a person building the table would:
# headers first, then the data rows, with the "Id" header added in
# after a couple of rows were already present.
$iv.Range("B1").Value = "Item"
$iv.Range("C1").Value = "Value"

$iv.Range("A2").Value = 0
$iv.Range("B2").Value = "BrokenSoul"
$iv.Range("C2").Value = 1

$iv.Range("A3").Value = 1
$iv.Range("B3").Value = "WhiteSoul"
$iv.Range("C3").Value = 3

$iv.Range("A1").Value = "Id"

$iv.Range("A4").Value = 2
$iv.Range("B4").Value = "BlueVioletSoult"
$iv.Range("C4").Value = 7

$iv.Range("A5").Value = 3
$iv.Range("B5").Value = "VioletSoul"
$iv.Range("C5").Value = 18

$iv.Range("A6").Value = 4
$iv.Range("B6").Value = "BlueSoul"
$iv.Range("C6").Value = 37

$iv.Range("A7").Value = 5
$iv.Range("B7").Value = "RedSoul"
$iv.Range("C7").Value = 57

$iv.Columns("B:B").ColumnWidth = 15

# ---------------------------------------------------------------------
# 2) Crystals sheet: a header row (18) + 5 data rows (19-23) that pull
#    item name/value via INDEX/MATCH against ItemValue, then compute a
#    "RealValue" and the ratio to the previous row.
# ---------------------------------------------------------------------
$cr = $wb.Worksheets.Item("Crystals")

$cr.Range("A18").Value = "Id"
$cr.Range("B18").Value = "ItemName"
$cr.Range("C18").Value = "Value"
$cr.Range("D18").Value = "input 1"
$cr.Range("E18").Value = "input 1 name"
$cr.Range("F18").Value = "input 1 value"
$cr.Range("G18").Value = "input 1"
$cr.Range("H18").Value = "input 1 name"
$cr.Range("I18").Value = "input 1 value"

# Back to ItemValue for the last lookup row before finishing the header.
$iv.Range("A8").Value = 6
$iv.Range("B8").Value = "nothing"
$iv.Range("C8").Value = 0

$cr.Range("J18").Value = "RealValue"
$cr.Range("K18").Value = 1

$cr.Range("A18:K18").Font.Bold = $true

# Row 19
$cr.Range("A19").Value = 1
$cr.Range("B19").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!A19,ItemValue!A:A))"
$cr.Range("C19").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!A19,ItemValue!A:A))"
$cr.Range("D19").Value = 0
$cr.Range("E19").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!D19,ItemValue!A:A))"
$cr.Range("F19").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!D19,ItemValue!A:A))"
$cr.Range("G19").Value = 6
$cr.Range("H19").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!G19,ItemValue!A:A))"
$cr.Range("I19").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!G19,ItemValue!A:A))"
$cr.Range("J19").Formula = "=C19-F19-I19"
$cr.Range("K19").Formula = "=J19/K18"

# Row 20
$cr.Range("A20").Value = 2
$cr.Range("B20").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!A20,ItemValue!A:A))"
$cr.Range("C20").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!A20,ItemValue!A:A))"
$cr.Range("D20").Value = 1
$cr.Range("E20").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!D20,ItemValue!A:A))"
$cr.Range("F20").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!D20,ItemValue!A:A))"
$cr.Range("G20").Value = 6
$cr.Range("H20").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!G20,ItemValue!A:A))"
$cr.Range("I20").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!G20,ItemValue!A:A))"
$cr.Range("J20").Formula = "=C20-F20-I20"
$cr.Range("K20").Formula = "=J20/J19"

# Row 21
$cr.Range("A21").Value = 3
$cr.Range("B21").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!A21,ItemValue!A:A))"
$cr.Range("C21").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!A21,ItemValue!A:A))"
$cr.Range("D21").Value = 1
$cr.Range("E21").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!D21,ItemValue!A:A))"
$cr.Range("F21").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!D21,ItemValue!A:A))"
$cr.Range("G21").Value = 2
$cr.Range("H21").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!G21,ItemValue!A:A))"
$cr.Range("I21").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!G21,ItemValue!A:A))"
$cr.Range("J21").Formula = "=C21-F21-I21"
$cr.Range("K21").Formula = "=J21/J20"

# Row 22
$cr.Range("A22").Value = 4
$cr.Range("B22").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!A22,ItemValue!A:A))"
$cr.Range("C22").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!A22,ItemValue!A:A))"
$cr.Range("D22").Value = 1
$cr.Range("E22").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!D22,ItemValue!A:A))"
$cr.Range("F22").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!D22,ItemValue!A:A))"
$cr.Range("G22").Value = 3
$cr.Range("H22").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!G22,ItemValue!A:A))"
$cr.Range("I22").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!G22,ItemValue!A:A))"
$cr.Range("J22").Formula = "=C22-F22-I22"
$cr.Range("K22").Formula = "=J22/J21"

# Row 23
$cr.Range("A23").Value = 5
$cr.Range("B23").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!A23,ItemValue!A:A))"
$cr.Range("C23").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!A23,ItemValue!A:A))"
$cr.Range("D23").Value = 2
$cr.Range("E23").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!D23,ItemValue!A:A))"
$cr.Range("F23").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!D23,ItemValue!A:A))"
$cr.Range("G23").Value = 3
$cr.Range("H23").Formula = "=INDEX(ItemValue!B:B,MATCH(Crystals!G23,ItemValue!A:A))"
$cr.Range("I23").Formula = "=INDEX(ItemValue!C:C,MATCH(Crystals!G23,ItemValue!A:A))"
$cr.Range("J23").Formula = "=C23-F23-I23"
$cr.Range("K23").Formula = "=J23/J22"

$cr.Columns("D:D").ColumnWidth = 12
$cr.Columns("E:E").ColumnWidth = 12.5703125
$cr.Columns("F:F").ColumnWidth = 12.42578125
$cr.Columns("H:H").ColumnWidth = 12.5703125
$cr.Columns("I:I").ColumnWidth = 12.42578125

$cr.PageSetup.PaperSize = 9
$cr.PageSetup.Orientation = 1

$cr.Range("K23").Select()
$iv.Range("J12").Select()
$cr.Select()

Write-Output "done"
